# Weekly update: insert a new data row after row 2 (Fruta / hortaliza, semanal).
# This shifts the previous rows 3-15 down to rows 4-16, and the new row 3
# receives the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, pushing existing rows 3..15 down to 4..16.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the latest data point.
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 45272
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100103
$ws.Range("H3").Value = "Frutos de hueso (carozo)"
$ws.Range("I3").Value = 100103003
$ws.Range("J3").Value = "Damasco"
$ws.Range("K3").Value = "Castle Brite"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 300
$ws.Range("N3").Value = 27000
$ws.Range("O3").Value = 28000
$ws.Range("P3").Value = 27500
$ws.Range("Q3").Value = "`$/caja 18 kilos granel"
$ws.Range("R3").Value = "Región de Coquimbo"
$ws.Range("S3").Value = 1528
$ws.Range("T3").Value = 18
